$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2
$ws.Range("B3").Value = 4
$ws.Range("B4").Value = 6
$ws.Range("B5").Value = 8
$ws.Range("B6").Value = 10
$ws.Range("B7").Value = 12

$ws.Range("B8").Select()
